$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $val)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "26.559.61"
Set-TextValue $ws "E2" "  -0.34%  "
Set-TextValue $ws "D3" "1.634.43"
Set-TextValue $ws "E3" "  +0.18%  "
Set-TextValue $ws "E4" "  +0.09%  "
Set-TextValue $ws "D5" "213.59"
Set-TextValue $ws "E5" "  +0.54%  "
Set-TextValue $ws "D6" "0.504"
Set-TextValue $ws "E6" "  +1.92%  "
Set-TextValue $ws "E7" "  +0.19%  "
Set-TextValue $ws "E8" "  -0.49%  "
Set-TextValue $ws "D9" "0.0625"
Set-TextValue $ws "E9" "  +0.33%  "
Set-TextValue $ws "D10" "18.83"
Set-TextValue $ws "E10" "  -0.76%  "
Set-TextValue $ws "E11" "  +0.21%  "
Set-TextValue $ws "D12" "1.861.97"
Set-TextValue $ws "E12" "  +0.16%  "
Set-TextValue $ws "D13" "1.660.28"
Set-TextValue $ws "E13" "  +1.72%  "
Set-TextValue $ws "E14" "  +1.77%  "
Set-TextValue $ws "D15" "0.524"
Set-TextValue $ws "E15" "  -0.50%  "
Set-TextValue $ws "D16" "65.28"
Set-TextValue $ws "E16" "  +3.56%  "
Set-TextValue $ws "D17" "26.577.20"
Set-TextValue $ws "D18" "0.0₃0742"
Set-TextValue $ws "E18" "  +0.35%  "
Set-TextValue $ws "D19" "215.70"
Set-TextValue $ws "E19" "  +2.92%  "
Set-TextValue $ws "E20" "  +0.17%  "
Set-TextValue $ws "D22" "6.26"
Set-TextValue $ws "E22" "  +1.41%  "
Set-TextValue $ws "D23" "9.34"
Set-TextValue $ws "E23" "  -0.86%  "
Set-TextValue $ws "D24" "2.22"
Set-TextValue $ws "E24" "  +14.53%  "
Set-TextValue $ws "D25" "147.34"
Set-TextValue $ws "E25" "  +0.25%  "
Set-TextValue $ws "D26" "1.01"
Set-TextValue $ws "E26" "  +0.25%  "
Set-TextValue $ws "E27" "  -0.50%  "
Set-TextValue $ws "D28" "6.90"
Set-TextValue $ws "E28" "  +1.46%  "
Set-TextValue $ws "D29" "15.63"
Set-TextValue $ws "E30" "  -1.49%  "
Set-TextValue $ws "E31" "  -0.46%  "
Set-TextValue $ws "E32" "  +3.89%  "
Set-TextValue $ws "E33" "  +0.93%  "
Set-TextValue $ws "D34" "1.260.44"
Set-TextValue $ws "E34" "  +8.00%  "
Set-TextValue $ws "E36" "  +0.15%  "
Set-TextValue $ws "E37" "  +4.30%  "
Set-TextValue $ws "E38" "  +1.39%  "
Set-TextValue $ws "E39" "  +0.23%  "
Set-TextValue $ws "E40" "  -1.25%  "
Set-TextValue $ws "D41" "2.27"
Set-TextValue $ws "E41" "  -1.96%  "
Set-TextValue $ws "E42" "  +0.73%  "
Set-TextValue $ws "E43" "  -0.44%  "
Set-TextValue $ws "D44" "1.771.07"
Set-TextValue $ws "E44" "  +0.07%  "
Set-TextValue $ws "D45" "93.41"
Set-TextValue $ws "E45" "  +1.46%  "
Set-TextValue $ws "E46" "  +3.16%  "
Set-TextValue $ws "D47" "54.99"
Set-TextValue $ws "E47" "  +0.61%  "
Set-TextValue $ws "D48" "0.0₆0103"
Set-TextValue $ws "E48" "  -1.73%  "
Set-TextValue $ws "D49" "0.0512"
Set-TextValue $ws "E49" "  +0.31%  "
Set-TextValue $ws "E50" "  -0.02%  "
Set-TextValue $ws "E51" "  -0.42%  "
